$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17, column A: was text "002", becomes the number 2 (time-to-go-off-road
# prediction work correcting the experimenter id to numeric).
$ws.Cells.Item(17, 1).Value = 2

# New rows 18-37: experimenter_id (numeric 2), condition, math_score,
# num_collisions, num_offroad_events, completion_time
$newRows = @(
    @(2, "Visual Only", 0, 0, 0, 24.55),
    @(2, "Visual Only", 0, 0, 5, 24.78),
    @(2, "Visual Only", 0, 0, 0, 10.36),
    @(2, "Visual Only", 0, 0, 10, 82.81999999999999),
    @(2, "Visual Only", 0, 0, 4, 29.01),
    @(2, "Visual Only", 0, 0, 0, 2.72),
    @(2, "Visual Only", 0, 0, 0, 7.39),
    @(2, "Visual Only", 0, 0, 0, 2.61),
    @(2, "Visual Only", 0, 0, 0, 2.9),
    @(2, "Visual Only", 0, 0, 8, 66.64),
    @(2, "No Feedback", 0, 0, 7, 26.32),
    @(2, "No Feedback", 0, 0, 11, 33.44),
    @(2, "No Feedback", 0, 0, 4, 48.81),
    @(2, "No Feedback", 0, 0, 4, 26.26),
    @(2, "No Feedback", 0, 0, 5, 48.3),
    @(2, "No Feedback", 0, 0, 12, 44.22),
    @(2, "No Feedback", 0, 0, 1, 23.18),
    @(2, "No Feedback", 0, 0, 0, 5.78),
    @(2, "No Feedback", 0, 0, 2, 12.9),
    @(2, "No Feedback", 0, 0, 7, 34.27)
)

$r = 18
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}

# Final row 38 keeps the experimenter id as the zero-padded text "002", like
# the original row 17 used to before it was corrected. A leading apostrophe
# forces Excel to store it as text rather than coercing it to the number 2.
$ws.Cells.Item(38, 1).Value = "'002"
$ws.Cells.Item(38, 2).Value = "No Feedback"
$ws.Cells.Item(38, 3).Value = 0
$ws.Cells.Item(38, 4).Value = 0
$ws.Cells.Item(38, 5).Value = 6
$ws.Cells.Item(38, 6).Value = 41.92
